$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected and needs to be inserted in its
# date-sorted position as row 289, pushing the existing rows 289-300 down
# to 290-301 (the sheet's used range grows from A1:R300 to A1:R301).
$ws.Rows("289:289").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A289").Value = 6
$ws.Range("B289").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C289").Value = "Metropolitana"
$ws.Range("D289").Value2 = 44509
$ws.Range("E289").Value = 13
$ws.Range("F289").Value = 100112043
$ws.Range("G289").Value = "Pepino ensalada"
$ws.Range("H289").Value = "Sin especificar"
$ws.Range("I289").Value = "Primera"
$ws.Range("J289").Value = 40
$ws.Range("K289").Value = 5000
$ws.Range("L289").Value = 6000
$ws.Range("M289").Value = 5575
$ws.Range("N289").Value = "`$/caja 70 unidades"
$ws.Range("O289").Value = "Provincia de Huasco"
$ws.Range("P289").Value = 80
$ws.Range("Q289").Value = 70
$ws.Range("R289").Value = "Hortaliza"
